$d = $word.ActiveDocument

$d.Content.Find.Execute("74÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "27÷6=", 2)
$d.Content.Find.Execute("42÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "97÷2=", 2)
$d.Content.Find.Execute("19÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "21÷3=", 2)
$d.Content.Find.Execute("93÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "12÷7=", 2)
$d.Content.Find.Execute("42÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "35÷3=", 2)
$d.Content.Find.Execute("34÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "12÷7=", 2)
$d.Content.Find.Execute("34÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "32÷5=", 2)
$d.Content.Find.Execute("29÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "92÷3=", 2)
$d.Content.Find.Execute("35÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "20÷8=", 2)
$d.Content.Find.Execute("14÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "46÷9=", 2)
$d.Content.Find.Execute("29÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "21÷6=", 2)
$d.Content.Find.Execute("16÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "34÷5=", 2)
$d.Content.Find.Execute("63÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "72÷6=", 2)
$d.Content.Find.Execute("94÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "70÷3=", 2)
$d.Content.Find.Execute("56÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "67÷8=", 2)
$d.Content.Find.Execute("14÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "81÷4=", 2)
$d.Content.Find.Execute("56÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "97÷6=", 2)
$d.Content.Find.Execute("18÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "62÷8=", 2)
$d.Content.Find.Execute("86÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "16÷3=", 2)
$d.Content.Find.Execute("60÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "12÷5=", 2)
$d.Content.Find.Execute("57÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "63÷9=", 2)
$d.Content.Find.Execute("31÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "18÷9=", 2)
$d.Content.Find.Execute("91÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "47÷7=", 2)
$d.Content.Find.Execute("30÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "18÷9=", 2)
$d.Content.Find.Execute("12÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "39÷8=", 2)
